# Insert a new data row at row 112 (pushes the existing rows 112-151 down to
# 113-152, and Excel auto-extends the used range / dimension accordingly),
# then populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(112).Insert()

$ws.Range("A112").Value = 10
$ws.Range("B112").Value = 'Vega Modelo de Temuco'
$ws.Range("C112").Value = 'La Araucanía'
$ws.Range("D112").Value = 44524
$ws.Range("E112").Value = 9
$ws.Range("F112").Value = 100112052
$ws.Range("G112").Value = 'Albahaca'
$ws.Range("H112").Value = 'Sin especificar'
$ws.Range("I112").Value = 'Primera'
$ws.Range("J112").Value = 30
$ws.Range("K112").Value = 3500
$ws.Range("L112").Value = 3500
$ws.Range("M112").Value = 3500
$ws.Range("N112").Value = '$/paquete'
$ws.Range("O112").Value = 'Región del Maule'
$ws.Range("P112").Value = 3500
$ws.Range("Q112").Value = 1
$ws.Range("R112").Value = 'Hortaliza'
